$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.103.56"
$ws.Range("E2").Value = "  -2.39%  "

$ws.Range("D3").Value = "'3.592.86"
$ws.Range("E3").Value = "  -2.73%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'623.36"
$ws.Range("E5").Value = "  -7.54%  "

$ws.Range("D6").Value = "'156.27"
$ws.Range("E6").Value = "  -3.19%  "

$ws.Range("D7").Value = "'3.588.52"
$ws.Range("E7").Value = "  -2.74%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'0.487"
$ws.Range("E9").Value = "  -2.54%  "

$ws.Range("D10").Value = "'0.142"
$ws.Range("E10").Value = "  -3.28%  "

$ws.Range("D11").Value = "'7.00"
$ws.Range("E11").Value = "  -1.60%  "

$ws.Range("D12").Value = "'0.433"
$ws.Range("E12").Value = "  -2.44%  "

$ws.Range("D13").Value = "'0.0000225"
$ws.Range("E13").Value = "  -4.05%  "

$ws.Range("D14").Value = "'4.195.92"
$ws.Range("E14").Value = "  -2.86%  "

$ws.Range("D15").Value = "'31.90"
$ws.Range("E15").Value = "  -4.29%  "

$ws.Range("D16").Value = "'3.604.12"
$ws.Range("E16").Value = "  -2.38%  "

$ws.Range("D17").Value = "'68.112.07"
$ws.Range("E17").Value = "  -2.31%  "

$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("D19").Value = "'6.42"
$ws.Range("E19").Value = "  -1.60%  "

$ws.Range("D20").Value = "'15.61"
$ws.Range("E20").Value = "  -3.72%  "

$ws.Range("D21").Value = "'9.90"
$ws.Range("E21").Value = "  +0.93%  "

$ws.Range("D22").Value = "'456.09"
$ws.Range("E22").Value = "  -3.50%  "

$ws.Range("D23").Value = "'0.641"
$ws.Range("E23").Value = "  -1.49%  "

$ws.Range("D24").Value = "'77.71"
$ws.Range("E24").Value = "  -3.04%  "

$ws.Range("D25").Value = "'3.733.71"
$ws.Range("E25").Value = "  -2.82%  "

$ws.Range("E26").Value = "  +0.06%  "

$ws.Range("D27").Value = "'10.73"
$ws.Range("E27").Value = "  -2.45%  "

$ws.Range("D28").Value = "'0.0000117"
$ws.Range("E28").Value = "  -9.05%  "

$ws.Range("D29").Value = "'8.42"
$ws.Range("E29").Value = "  -7.95%  "

$ws.Range("D30").Value = "'2.59"
$ws.Range("E30").Value = "  -4.23%  "

$ws.Range("D31").Value = "'1.63"
$ws.Range("E31").Value = "  -5.67%  "

$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").Value = "'26.07"
$ws.Range("E33").Value = "  -3.11%  "

$ws.Range("D34").Value = "'1.92"
$ws.Range("E34").Value = "  -5.70%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.159"
$ws.Range("E35").Value = "  -5.46%  "

$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "'3.590.10"
$ws.Range("E36").Value = "  -2.70%  "

$ws.Range("D37").Value = "'6.23"
$ws.Range("E37").Value = "  -4.76%  "

$ws.Range("D38").Value = "'8.16"
$ws.Range("E38").Value = "  -4.32%  "

$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("D41").Value = "'176.89"
$ws.Range("E41").Value = "  +0.58%  "

$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'5.64"
$ws.Range("E42").Value = "  -8.24%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'2.15"
$ws.Range("E43").Value = "  -5.85%  "

$ws.Range("D44").Value = "'0.0881"
$ws.Range("E44").Value = "  -3.34%  "

$ws.Range("D45").Value = "'0.903"
$ws.Range("E45").Value = "  -3.53%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'28.94"
$ws.Range("E46").Value = "  +3.42%  "

$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").Value = "'45.94"
$ws.Range("E47").Value = "  -2.23%  "

$ws.Range("D48").Value = "'2.59"
$ws.Range("E48").Value = "  -6.83%  "

$ws.Range("D49").Value = "'7.71"
$ws.Range("E49").Value = "  -2.59%  "

$ws.Range("B50").Value = "SuiNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D50").Value = "'1.02"
$ws.Range("E50").Value = "  -5.99%  "

$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Value = "'1.19"
$ws.Range("E51").Value = "  -7.89%  "
